# Update the "Field" value in row 2 for each Groups-related sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 (Country_Group_N): Field for the first condition changes
# from "Location Number" to "Country"
$ws1.Range("C2").Value = "Country"

# Sheet2 (Country_Group_A): Field for the first condition changes
# from "Phone" to "Postal Code"
$ws2.Range("C2").Value = "Postal Code"

# Update the selected cell on each sheet to match where the edits were made
$ws1.Range("C2").Select()
$ws2.Range("C4").Select()
